# Revert the earlier "edited excel files" commit: the router IP addresses
# in Sheet1 column F were given a "/24" CIDR suffix; put them back to the
# plain dotted-quad form. Also restore the previously-selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = "192.168.11.1"
$ws.Range("F3").Value = "192.168.11.2"
$ws.Range("F4").Value = "192.168.10.1"
$ws.Range("F5").Value = "10.1.1.1"

# Restore the active-cell selection that was in effect before the edit.
$ws.Range("F5").Select() | Out-Null
